# "add comments in file manager"
# The underlying data edit (as captured by the OOXML diff) is:
#  - vocabulary sheet:  trailing test rows 100-106 were removed (dimension A1:C106 -> A1:C99)
#  - categories sheet:  trailing test rows 9-12 were removed   (dimension A1:B12 -> A1:B8)
#  - test results sheet: F15 corrected from text "7.62" to the number 7.62,
#                         and a new test run was logged as row 16
#  - the "categories" sheet became the active/selected tab (was "vocabulary")

$wb = $excel.ActiveWorkbook
$wsVocabulary   = $wb.Worksheets.Item(1)   # "vocabulary"
$wsCategories   = $wb.Worksheets.Item(2)   # "categories"
$wsTestResults  = $wb.Worksheets.Item(3)   # "test results"

# --- vocabulary: drop the leftover scratch/test entries (rows 100-106) -------
$wsVocabulary.Rows("100:106").Delete()

# --- categories: drop the leftover scratch/test entries (rows 9-12) ----------
$wsCategories.Rows("9:12").Delete()

# --- test results: fix F15 (was stored as text "7.62", should be numeric) ----
$wsTestResults.Range("F15").Value = 7.62

# --- test results: append the new test run as row 16 -------------------------
$wsTestResults.Range("A16").Value = 15
$wsTestResults.Range("B16").Value = "filip"
$wsTestResults.Range("C16").Value = "02-02-2025 13:08:52"
$wsTestResults.Range("D16").Value = "PL->EN"
$wsTestResults.Range("E16").Value = 60

# F16/G16/H16/I16 mirror the source data's quirky text-typed numbers —
# force plain-text storage (no residual number format) so the written value
# matches the literal string rather than being auto-coerced to a number.
$wsTestResults.Range("F16").NumberFormat = "@"
$wsTestResults.Range("F16").Value = "127.13"
$wsTestResults.Range("F16").ClearFormats()

$wsTestResults.Range("G16").Value = "(6/30)"

$wsTestResults.Range("H16").NumberFormat = "@"
$wsTestResults.Range("H16").Value = "20.00%"
$wsTestResults.Range("H16").ClearFormats()

$wsTestResults.Range("I16").Value = "Games Remastering, Processors, Graphic cards"

# --- categories becomes the active/selected sheet tab -------------------------
$wsCategories.Activate()
$wsCategories.Range("A9:B12").Select()
